$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove the old "LEONARDO FELIX SEÑAS DIAZ" worker rows (rows 16-21).
#    After the delete, the remaining "ELIAS DAVID RIPOLL POLO" rows (old
#    22-130, 109 rows) shift up to become rows 16-124, and the closing
#    (bottom-border) row style that used to live on row 130 now lives on
#    row 124.
# ---------------------------------------------------------------------------
$ws.Rows("16:21").Delete()

# ---------------------------------------------------------------------------
# 2. We need one extra data row for the newly added period (2508), so the
#    data block grows from 109 to 110 rows (16-125). Insert a fresh row
#    right below the current last data row (124).
# ---------------------------------------------------------------------------
$ws.Rows("125:125").Insert()

# Copy the "closing" (bottom-border) formatting that is currently on row 124
# down onto the brand new row 125, then restore row 124 back to the regular
# "interior" row formatting (copied from row 123).
$ws.Range("B124:J124").Copy()
$ws.Range("B125:J125").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("B123:J123").Copy()
$ws.Range("B124:J124").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3. Fill the 110 data rows (16-125) with the single remaining worker's
#    records, sorted in ascending period order (1607 .. 2508), including the
#    brand-new 2508 period.
# ---------------------------------------------------------------------------
$periods = "1607","1608","1609","1610","1611","1612","1701","1702","1703","1704","1705","1706","1707","1708","1709","1710","1711","1712","1801","1802","1803","1804","1805","1806","1807","1808","1809","1810","1811","1812","1901","1902","1903","1904","1905","1906","1907","1908","1909","1910","1911","1912","2001","2002","2003","2004","2005","2006","2007","2008","2009","2010","2011","2012","2101","2102","2103","2104","2105","2106","2107","2108","2109","2110","2111","2112","2201","2202","2203","2204","2205","2206","2207","2208","2209","2210","2211","2212","2301","2302","2303","2304","2305","2306","2307","2308","2309","2310","2311","2312","2401","2402","2403","2404","2405","2406","2407","2408","2409","2410","2411","2412","2501","2502","2503","2504","2505","2506","2507","2508"

$startRow = 16
for ($i = 0; $i -lt $periods.Length; $i++) {
    $r = $startRow + $i
    $ws.Range("B" + $r).Value = "CC"
    $ws.Range("C" + $r).Value = "73205067"
    $ws.Range("D" + $r).Value = "ELIAS DAVID RIPOLL POLO"
    $ws.Range("E" + $r).Value = $periods[$i]
    $ws.Range("F" + $r).Value = 60000
    $ws.Range("G" + $r).Value = 1500000
}

# ---------------------------------------------------------------------------
# 4. Update the summary figures at the top of the sheet:
#      - total overdue amount (VALOR MORA)
#      - number of workers
#      - number of overdue periods
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = 6600000
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 110

# ---------------------------------------------------------------------------
# 5. Re-fit column D (worker name) now that the long "LEONARDO FELIX SEÑAS
#    DIAZ" name is gone and every row uses the shorter "ELIAS DAVID RIPOLL
#    POLO" name.
# ---------------------------------------------------------------------------
$ws.Columns("D").AutoFit()
